$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers
$ws.Range("J1").Value = "(*)Commodity Category"
$ws.Range("K1").Value = "Sub Category"

# Add new headers for columns W..AC (23..29)
# Shared-string table order follows entry order: Length, Width, Initail Quantity, Feet, Meters, Measurement, As Of Date
$ws.Cells.Item(1, 23).Value = "Length"
$ws.Cells.Item(1, 24).Value = "Width"
$ws.Cells.Item(1, 28).Value = "Initail Quantity "
$ws.Cells.Item(1, 25).Value = "Feet"
$ws.Cells.Item(1, 26).Value = "Meters"
$ws.Cells.Item(1, 27).Value = "Measurement"
$ws.Cells.Item(1, 29).Value = "As Of Date"

# Copy the style of the last existing header (V1) to the new header cells
$srcRange = $ws.Range("V1")
$dstRange = $ws.Range("W1:AC1")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122) # xlPasteFormats

$ws.Range("AC1").Select()
